$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.920.34'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '3.508.65'
$ws.Range('E3').Value = '  -0.50%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '600.47'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.48%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '196.23'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +6.83%  '
$ws.Range('E7').Value = '  +1.88%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('E10').Value = '  +2.22%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '54.05'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('E12').Value = '  -1.94%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '9.54'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('D14').Value = '4.061.46'
$ws.Range('E14').Value = '  -0.62%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '607.17'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +3.16%  '
$ws.Range('D16').Value = '70.070.95'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('E17').Value = '  +0.89%  '
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').Value = '3.511.33'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('E20').Value = '  +0.84%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.994'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.74%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '17.91'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +2.68%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '104.23'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +7.72%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.15'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +6.53%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.59'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.89%  '
$ws.Range('E26').Value = '  +3.55%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.98'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.61%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '9.69'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +1.68%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '33.67'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +5.39%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.58'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +28.12%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.10'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.92%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '12.67'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.77%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.116'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.23%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '63.20'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('D35').Value = '0.0₃0827'
$ws.Range('E35').Value = '  +6.71%  '
$ws.Range('D36').Value = '3.737.00'
$ws.Range('E36').Value = '  +5.70%  '
$ws.Range('E37').Value = '  -4.83%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('E39').Value = '  -1.59%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '36.80'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('E41').Value = '  +1.90%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '499.43'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -6.33%  '
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.83'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -2.75%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.140'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('E48').Value = '  +0.33%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '8.70'
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.000245'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.01%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '130.88'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.36%  '
